$wb = $excel.ActiveWorkbook
$mainData = $wb.Worksheets.Item("Main Data")

# ---------------------------------------------------------------------------
# 1. Fill in problems 27-40 (rows 30-43) on "Main Data": column A = SN number,
#    column B = "Problem<N>" name.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 14; $i++) {
    $row = 30 + $i
    $num = 27 + $i
    $mainData.Cells.Item($row, 1).Value = $num
    $mainData.Cells.Item($row, 2).Value = "Problem$num"
}

# ---------------------------------------------------------------------------
# 2. Adjust the "Main Data" sheet view: zoom + selection.
# ---------------------------------------------------------------------------
$mainData.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 55
$mainData.Columns("E:E").Select()

# ---------------------------------------------------------------------------
# 3. Add the new "ChatGPT" worksheet right after "Main Data".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $mainData)
$newSheet.Name = "ChatGPT"

$newSheet.Range("A1").Value = "Problem Name"
$newSheet.Range("B1").Value = "Observation"
$newSheet.Range("A2").Value = "Word Search II"
$newSheet.Range("E2").Value = "This page has the list of probelems that ChatGPT could not solve."
$newSheet.Range("E3").Value = "Maybe these can be a measure for AGI"

# ---------------------------------------------------------------------------
# 4. Turn A1:B2 into "Table1" styled with TableStyleLight1.
# ---------------------------------------------------------------------------
$tbl = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:B2"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight1"

# ---------------------------------------------------------------------------
# 5. Make "ChatGPT" the active/selected tab (matches activeTab=1 in the diff).
# ---------------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("E4").Select()
